$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.800.63"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").Value = "2.077.73"
$ws.Range("E3").Value = "  -0.46%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "233.10"
$ws.Range("E5").Value = "  -0.87%  "

$ws.Range("D6").Value = "0.625"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").Value = "58.64"
$ws.Range("E7").Value = "  -1.42%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").Value = "0.0786"
$ws.Range("E10").Value = "  -0.69%  "

$ws.Range("E11").Value = "  +3.37%  "

$ws.Range("D12").Value = "2.384.23"

$ws.Range("D13").Value = "14.79"
$ws.Range("E13").Value = "  +0.32%  "

$ws.Range("D14").Value = "21.13"
$ws.Range("E14").Value = "  -1.59%  "

$ws.Range("D15").Value = "0.781"
$ws.Range("E15").Value = "  +0.98%  "

$ws.Range("E16").Value = "  +0.78%  "

$ws.Range("D17").Value = "2.070.91"
$ws.Range("E17").Value = "  -0.83%  "

$ws.Range("D18").Value = "37.688.49"
$ws.Range("E18").Value = "  -0.29%  "

$ws.Range("D19").Value = "6.13"
$ws.Range("E19").Value = "  -1.81%  "

$ws.Range("D20").Value = "71.64"
$ws.Range("E20").Value = "  -0.32%  "

$ws.Range("E21").Value = "  +1.47%  "

$ws.Range("D22").Value = "229.47"
$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").Value = "2.38"
$ws.Range("E24").Value = "  -1.37%  "

$ws.Range("E25").Value = "  -0.30%  "

$ws.Range("E26").Value = "  +7.30%  "

$ws.Range("D27").Value = "171.76"
$ws.Range("E27").Value = "  +0.47%  "

$ws.Range("E28").Value = "  -0.43%  "

$ws.Range("E29").Value = "  -2.17%  "

$ws.Range("D30").Value = "19.43"
$ws.Range("E30").Value = "  -0.93%  "

$ws.Range("D31").Value = "0.121"
$ws.Range("E31").Value = "  +1.03%  "

$ws.Range("E32").Value = "  +0.32%  "

$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("D34").Value = "4.68"
$ws.Range("E34").Value = "  -0.54%  "

$ws.Range("D35").Value = "2.45"
$ws.Range("E35").Value = "  -2.67%  "

$ws.Range("E36").Value = "  -0.13%  "

$ws.Range("D37").Value = "3.41"
$ws.Range("E37").Value = "  -3.16%  "

$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("D39").Value = "5.42"
$ws.Range("E39").Value = "  -0.96%  "

$ws.Range("D40").Value = "0.0233"
$ws.Range("E40").Value = "  +7.92%  "

$ws.Range("D41").Value = "100.99"
$ws.Range("E41").Value = "  +1.28%  "

$ws.Range("D42").Value = "0.0977"
$ws.Range("E42").Value = "  -1.02%  "

$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").Value = "2.93"
$ws.Range("E43").Value = "  -0.65%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "17.06"
$ws.Range("E44").Value = "  +5.22%  "

$ws.Range("D45").Value = "1.451.56"
$ws.Range("E45").Value = "  -0.93%  "

$ws.Range("E46").Value = "  -1.83%  "

$ws.Range("E47").Value = "  -0.71%  "

$ws.Range("E48").Value = "  -5.27%  "

$ws.Range("D49").Value = "7.41"
$ws.Range("E49").Value = "  -0.95%  "

$ws.Range("E50").Value = "  -1.51%  "

$ws.Range("D51").Value = "2.269.69"
$ws.Range("E51").Value = "  -0.42%  "
